$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.993.93"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "2.048.10"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'248.56"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").Value = "'0.662"
$ws.Range("E6").Value = "  +1.82%  "
$ws.Range("D7").Value = "'57.77"
$ws.Range("E7").Value = "  +4.82%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.381"
$ws.Range("E9").Value = "  +1.92%  "
$ws.Range("D10").Value = "'0.0777"
$ws.Range("E10").Value = "  +2.23%  "
$ws.Range("D11").Value = "'0.107"
$ws.Range("E11").Value = "  +2.02%  "
$ws.Range("D12").Value = "'15.76"
$ws.Range("E12").Value = "  +5.83%  "
$ws.Range("D13").Value = "2.351.15"
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("D14").Value = "'0.799"
$ws.Range("E14").Value = "  -1.25%  "
$ws.Range("D15").Value = "'5.58"
$ws.Range("E15").Value = "  +8.18%  "
$ws.Range("D16").Value = "2.048.13"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").Value = "37.036.29"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D19").Value = "'74.60"
$ws.Range("E19").Value = "  +4.12%  "
$ws.Range("D20").Value = "0.0₃0898"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "'5.33"
$ws.Range("E21").Value = "  +2.11%  "
$ws.Range("D22").Value = "'235.78"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'2.36"
$ws.Range("E24").Value = "  -1.26%  "
$ws.Range("D25").Value = "'2.20"
$ws.Range("E25").Value = "  +11.50%  "
$ws.Range("D26").Value = "'167.66"
$ws.Range("E26").Value = "  -0.79%  "
$ws.Range("D27").Value = "'9.12"
$ws.Range("E27").Value = "  +1.24%  "
$ws.Range("D28").Value = "'19.74"
$ws.Range("E28").Value = "  -0.97%  "
$ws.Range("E29").Value = "  +1.59%  "
$ws.Range("E30").Value = "  +6.95%  "
$ws.Range("D31").Value = "'4.67"
$ws.Range("E31").Value = "  +3.14%  "
$ws.Range("D32").Value = "'0.0612"
$ws.Range("E32").Value = "  -0.63%  "
$ws.Range("D33").Value = "'4.47"
$ws.Range("E33").Value = "  +3.95%  "
$ws.Range("D34").Value = "'0.0882"
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").Value = "'2.22"
$ws.Range("E36").Value = "  -1.02%  "
$ws.Range("E37").Value = "  -1.85%  "
$ws.Range("D38").Value = "'0.108"
$ws.Range("E38").Value = "  +5.90%  "
$ws.Range("E39").Value = "  +0.40%  "
$ws.Range("D40").Value = "'3.16"
$ws.Range("E40").Value = "  +13.05%  "
$ws.Range("D41").Value = "'4.93"
$ws.Range("E41").Value = "  +25.06%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0220"
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "'17.30"
$ws.Range("E43").Value = "  -3.87%  "
$ws.Range("E44").Value = "  -1.09%  "
$ws.Range("D45").Value = "'95.42"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D46").Value = "'2.44"
$ws.Range("E46").Value = "  +5.26%  "
$ws.Range("D47").Value = "1.279.33"
$ws.Range("E47").Value = "  -0.82%  "
$ws.Range("D48").Value = "'2.85"
$ws.Range("E48").Value = "  -1.61%  "
$ws.Range("D49").Value = "2.241.29"
$ws.Range("E49").Value = "  +0.72%  "
$ws.Range("D50").Value = "'6.67"
$ws.Range("E50").Value = "  -0.89%  "
$ws.Range("D51").Value = "'3.59"
$ws.Range("E51").Value = "  -10.98%  "
